$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "46.035.58"
Set-TextValue "E2" "  -0.64%  "
Set-TextValue "D3" "2.598.06"
Set-TextValue "E3" "  -0.30%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.17%  "
Set-TextValue "D5" "311.89"
Set-TextValue "E5" "  +1.88%  "
Set-TextValue "D6" "98.73"
Set-TextValue "E6" "  -1.75%  "
Set-TextValue "D7" "0.600"
Set-TextValue "E7" "  -0.47%  "
Set-TextValue "E8" "  -0.02%  "
Set-TextValue "D9" "0.583"
Set-TextValue "E9" "  +1.63%  "
Set-TextValue "D10" "39.04"
Set-TextValue "E10" "  +0.34%  "
Set-TextValue "D11" "54.16"
Set-TextValue "E11" "  -1.76%  "
Set-TextValue "D12" "0.0841"
Set-TextValue "E12" "  +0.60%  "
Set-TextValue "D13" "8.13"
Set-TextValue "E13" "  -0.22%  "
Set-TextValue "D14" "2.989.73"
Set-TextValue "E14" "  -0.85%  "
Set-TextValue "E15" "  +1.37%  "
Set-TextValue "D16" "2.594.36"
Set-TextValue "E16" "  -1.27%  "
Set-TextValue "D17" "0.918"
Set-TextValue "E17" "  +2.07%  "
Set-TextValue "D18" "14.85"
Set-TextValue "E18" "  +0.02%  "
Set-TextValue "D19" "46.077.00"
Set-TextValue "E19" "  -0.85%  "
Set-TextValue "D20" "0.0000102"
Set-TextValue "E20" "  +1.58%  "
Set-TextValue "D21" "6.73"
Set-TextValue "E21" "  +0.68%  "
Set-TextValue "D22" "12.80"
Set-TextValue "E22" "  -3.39%  "
Set-TextValue "D23" "296.92"
Set-TextValue "E23" "  +16.14%  "
Set-TextValue "D24" "72.99"
Set-TextValue "E24" "  +2.67%  "
Set-TextValue "D25" "3.05"
Set-TextValue "E25" "  +1.42%  "
Set-TextValue "D26" "2.23"
Set-TextValue "E26" "  +0.26%  "
Set-TextValue "D27" "29.60"
Set-TextValue "E27" "  +5.36%  "
Set-TextValue "E28" "  -0.03%  "
Set-TextValue "D29" "4.05"
Set-TextValue "E29" "  +0.77%  "
Set-TextValue "D30" "10.79"
Set-TextValue "E30" "  +3.36%  "
Set-TextValue "D31" "38.60"
Set-TextValue "E31" "  -3.26%  "
Set-TextValue "E32" "  -2.55%  "
Set-TextValue "D33" "6.23"
Set-TextValue "E33" "  +1.35%  "
Set-TextValue "D34" "3.57"
Set-TextValue "E34" "  -4.12%  "
Set-TextValue "D35" "155.76"
Set-TextValue "E35" "  +3.00%  "
Set-TextValue "D36" "0.0836"
Set-TextValue "E36" "  +0.58%  "
Set-TextValue "D37" "2.20"
Set-TextValue "E37" "  -5.49%  "
Set-TextValue "D38" "2.79"
Set-TextValue "E38" "  -5.64%  "
Set-TextValue "E39" "  +3.46%  "
Set-TextValue "E40" "  +1.20%  "
Set-TextValue "D41" "15.73"
Set-TextValue "E41" "  +0.66%  "
Set-TextValue "D42" "0.0331"
Set-TextValue "E42" "  +2.80%  "
Set-TextValue "D43" "3.58"
Set-TextValue "E43" "  -0.64%  "
Set-TextValue "D44" "21.28"
Set-TextValue "E44" "  +9.77%  "
Set-TextValue "D45" "3.95"
Set-TextValue "E45" "  -5.10%  "
Set-TextValue "D46" "2.108.40"
Set-TextValue "E46" "  +2.66%  "
Set-TextValue "D47" "98.64"
Set-TextValue "E47" "  +8.47%  "
Set-TextValue "E48" "  -0.05%  "
Set-TextValue "D49" "9.58"
Set-TextValue "E49" "  +4.32%  "
Set-TextValue "D50" "0.202"
Set-TextValue "E50" "  +1.30%  "
Set-TextValue "D51" "108.61"
Set-TextValue "E51" "  -1.00%  "
